# Currency changes between dates sheet implementation
#
# Fix the typo in the "Símbolo" column header (was "Símbulo") and
# update the active selection to D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "Símbulo" -> "Símbolo" in the header cell D1
$ws.Range("D1").Value = "Símbolo"

# Update the active cell/selection to D5
$ws.Range("D5").Select()
